# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interested-count) figures pulled from bilibili show pages.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 5515
$ws.Range("F4").Value = 637
$ws.Range("F6").Value = 830
$ws.Range("F7").Value = 31
$ws.Range("F8").Value = 358

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 45
$ws.Range("F3").Value = 17

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5515
$ws.Range("F4").Value = 637
$ws.Range("F6").Value = 830
$ws.Range("F7").Value = 31
$ws.Range("F8").Value = 45
$ws.Range("F9").Value = 358
$ws.Range("F13").Value = 17
